$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.139.42"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "3.542.06"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.49%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "

$ws.Range("D7").Value = "3.539.53"
$ws.Range("E7").Value = "  -0.43%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.15%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.501"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.16%  "

$ws.Range("E10").Value = "  -1.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.409"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.04%  "

$ws.Range("D13").Value = "4.144.22"
$ws.Range("E13").Value = "  -0.26%  "

$ws.Range("E14").Value = "  -2.45%  "

$ws.Range("E15").Value = "  -3.63%  "

$ws.Range("D16").Value = "3.553.61"
$ws.Range("E16").Value = "  +0.07%  "

$ws.Range("E17").Value = "  +1.55%  "

$ws.Range("D18").Value = "66.121.62"
$ws.Range("E18").Value = "  -0.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.63"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "419.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.61%  "

$ws.Range("E23").Value = "  -1.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.77"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.20%  "

$ws.Range("D25").Value = "3.682.79"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  +0.10%  "

$ws.Range("E27").Value = "  -2.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.11"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "

$ws.Range("E29").Value = "  -1.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.82"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("E31").Value = "  +0.03%  "

$ws.Range("D32").Value = "3.541.15"
$ws.Range("E32").Value = "  -0.17%  "

$ws.Range("E33").Value = "  +1.00%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.44"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.91%  "

$ws.Range("E36").Value = "  -3.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.88%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.29"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.29%  "

$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "174.10"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "

$ws.Range("B40").Value = "ImmutableX"
$ws.Range("C40").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.60"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.31%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0825"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.76%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.60%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.860"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.65"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "

$ws.Range("E45").Value = "  -5.88%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.01%  "

$ws.Range("E47").Value = "  -3.01%  "

$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("B49").Value = "ONDO"
$ws.Range("C49").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.90%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.33%  "
